{"js": "// Office.js (Word JavaScript API) script\n// Applies the edits described by the diff:\n//  1. Prepend \"Anexo 3: \" (not underlined) to the title paragraph, and relocate the\n//     \"_GoBack\" bookmark so it sits right after the new \"Anexo 3: \" text (before \"Gu\u00eda\").\n//  2. Rewrite four consecutive bullet questions with the text that now appears in the\n//     revised document (the paragraphs were effectively re-ordered/re-worded).\n//  3. Remove the trailing \"Conclusi\u00f3n\" bullet paragraph and the blank paragraph that\n//     immediately followed it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1. Title paragraph: \"Anexo 3: \" prefix + bookmark relocation ---------------\nconst titlePara = paragraphs.items[0];\nconst titleStart = titlePara.getRange(\"Start\");\ntitleStart.insertText(\"Anexo 3: \", Word.InsertLocation.before);\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark from the end of \"Segunda entrevista\" to right after\n// the freshly inserted \"Anexo 3: \" text (i.e. immediately before \"Gu\u00eda\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst anexoSearch = body.search(\"Anexo 3: \", { matchCase: true });\nanexoSearch.load(\"items\");\nawait context.sync();\nconst anexoEnd = anexoSearch.items[0].getRange(\"End\");\nanexoEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2. Rewrite the four bullet questions --------------------------------------\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst newTexts = {\n  12: \"\u00bfTiene en mente alguna alternativa sobre como le gustar\u00eda que fuera la distribuci\u00f3n, visualizaci\u00f3n y/o flujo del sistema?\",\n  13: \"La informaci\u00f3n que le brinde el sistema, \u00bfQu\u00e9 tan \u00fatil le parecer\u00eda compartirla? \u00bfCon quien compartir\u00eda este informaci\u00f3n?\",\n  14: \"\u00bfTiene alg\u00fan otro comentario o sugerencia sobre los bocetos, flujo y/o sistema?\",\n  15: \"\u00bfCu\u00e1l es su opini\u00f3n general del sistema?\"\n};\n\nfor (const idx of Object.keys(newTexts)) {\n  const p = paragraphs.items[Number(idx)];\n  const whole = p.getRange(\"Whole\");\n  whole.insertText(newTexts[idx], Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 3. Remove the \"Conclusi\u00f3n\" paragraph and the blank line after it ----------\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst conclusionIndex = paragraphs.items.findIndex((p) => p.text.trim() === \"Conclusi\u00f3n\");\nif (conclusionIndex !== -1) {\n  const conclusionPara = paragraphs.items[conclusionIndex];\n  const blankAfterPara = paragraphs.items[conclusionIndex + 1];\n\n  // Delete the blank paragraph first so the \"Conclusi\u00f3n\" paragraph's index/range\n  // stays valid for the second delete call.\n  if (blankAfterPara) {\n    blankAfterPara.delete();\n  }\n  conclusionPara.delete();\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Applies the edits described by the diff:\n#  1. Prepend \"Anexo 3: \" (not underlined) to the title paragraph, and relocate the\n#     \"_GoBack\" bookmark so it sits right after the new \"Anexo 3: \" text (before \"Gu\u00eda\").\n#  2. Rewrite four consecutive bullet questions with the text that now appears in the\n#     revised document (the paragraphs were effectively re-ordered/re-worded).\n#  3. Remove the trailing \"Conclusi\u00f3n\" bullet paragraph and the blank paragraph that\n#     immediately followed it.\n\n$d = $word.ActiveDocument\n\n# --- 1. Title paragraph: \"Anexo 3: \" prefix + bookmark relocation -----------------\n$titlePara = $d.Paragraphs(1)\n$titleStart = $d.Range($titlePara.Range.Start, $titlePara.Range.Start)\n$titleStart.InsertBefore(\"Anexo 3: \")\n\n# Move the \"_GoBack\" bookmark from the end of \"Segunda entrevista\" to right after\n# the freshly inserted \"Anexo 3: \" text (i.e. immediately before \"Gu\u00eda\").\n$d.Bookmarks(\"_GoBack\").Delete()\n\n$titlePara = $d.Paragraphs(1)\n$bookmarkPos = $titlePara.Range.Start + 9   # length of \"Anexo 3: \"\n$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n# --- 2. Rewrite the four bullet questions -----------------------------------------\nfunction Set-ParagraphText($paraIndex, $newText) {\n    $p = $d.Paragraphs($paraIndex)\n    $r = $p.Range\n    $r.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark\n    $r.Text = $newText\n}\n\nSet-ParagraphText 13 \"\u00bfTiene en mente alguna alternativa sobre como le gustar\u00eda que fuera la distribuci\u00f3n, visualizaci\u00f3n y/o flujo del sistema?\"\nSet-ParagraphText 14 \"La informaci\u00f3n que le brinde el sistema, \u00bfQu\u00e9 tan \u00fatil le parecer\u00eda compartirla? \u00bfCon quien compartir\u00eda este informaci\u00f3n?\"\nSet-ParagraphText 15 \"\u00bfTiene alg\u00fan otro comentario o sugerencia sobre los bocetos, flujo y/o sistema?\"\nSet-ParagraphText 16 \"\u00bfCu\u00e1l es su opini\u00f3n general del sistema?\"\n\n# --- 3. Remove the \"Conclusi\u00f3n\" paragraph and the blank line after it ------------\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $txt = $d.Paragraphs($i).Range.Text.Trim()\n    if ($txt -eq \"Conclusi\u00f3n\") {\n        $conclusionPara = $d.Paragraphs($i)\n        $nextPara = $d.Paragraphs($i + 1)\n        $deleteRange = $d.Range($conclusionPara.Range.Start, $nextPara.Range.End)\n        $deleteRange.Delete()\n        break\n    }\n}\n"}
